$wb = $excel.ActiveWorkbook
$errSheet = $wb.Worksheets.Item("Error Message Catalog")
$newSheet = $wb.Worksheets.Add($errSheet)

$newSheet.Range("G2").Value = "title"
$newSheet.Range("H2").Value = "titulo de la discusión"
$newSheet.Range("B3").Value = "comments"
$newSheet.Range("F5:F10").Merge()
$newSheet.Range("F5").Value = "id"
$newSheet.Range("F5:F10").HorizontalAlignment = -4108
$newSheet.Range("F5:F10").VerticalAlignment = -4108

Write-Host "done"
